$d = $word.ActiveDocument

# --- Edit 1: "Jan 2013 - Present" -> "Jan 2013 - " + "Dec 2018" (new run) ---
$r1 = $d.Content
$found1 = $r1.Find.Execute("Jan 2013 " + [char]0x2013 + " Present", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Jan 2013 - Present' text"
}
$prefix1 = "Jan 2013 " + [char]0x2013 + " "
$splitPos1 = $r1.Start + $prefix1.Length
$splitPoint1 = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("TempSplitDate", $splitPoint1) | Out-Null
$tailRange1 = $d.Range($splitPos1, $r1.End)
$tailRange1.Text = "Dec 2018"
$d.Bookmarks("TempSplitDate").Delete()

# --- Edit 2: "...Subject Matter Expert for RMX - a p" ->
#     "...Subject Matter Expert for" + _GoBack bookmark + " RMX - a p" ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("Subject Matter Expert for RMX " + [char]0x2013 + " a p", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'Subject Matter Expert for RMX - a p' text"
}
$marker = "for RMX"
$forIdx = $r2.Text.IndexOf($marker)
$splitPos2 = $r2.Start + $forIdx + 3
$splitPoint2 = $d.Range($splitPos2, $splitPos2)

# Moving/re-adding a bookmark named "_GoBack" relocates Word's special
# "last edit" bookmark from wherever it previously was to this new spot.
$d.Bookmarks.Add("_GoBack", $splitPoint2) | Out-Null
